# repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows to reflect repulled data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -1
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F34").Value = -2
$ws.Range("F38").Value = -9
$ws.Range("F41").Value = -2
$ws.Range("F43").Value = 4
$ws.Range("F45").Value = -3
